# Extra-Wide Flanged Button Head Screws - restructure header rows
#
# Summary of the change being applied:
#  - Two new rows are inserted above the existing data rows (which shift
#    down from rows 2-19 to rows 4-21; their content is untouched).
#  - Row 1 (previously text column headers) becomes a row of plain
#    numbers 0..10.
#  - A new row 2 is mostly blank, except E2 = "Drive".
#  - A new row 3 holds the text labels that used to live in row 1
#    (Lg., Threading, FlangeDia., Head Ht., Style, Size,
#    TensileStrength/psi, <blank>, Each) - but without the old
#    J1 ("thread_size") / K1 ("material_surface") values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original row-1 header text before we overwrite anything.
$oldHeaders = @{}
$cols = @("A","B","C","D","E","F","G","H","I","J","K")
foreach ($col in $cols) {
    $oldHeaders[$col] = $ws.Range($col + "1").Value()
}

# Insert two blank rows at the top of the data (pushes old row 2.. down to 4..)
$ws.Range("A2:A3").EntireRow.Insert()

# The newly inserted rows pick up the bold/boxed header formatting from
# row 1 above - strip that back to the plain/default look used by the
# rest of the data rows.
$ws.Range("A2:K3").ClearFormats()

# Row 1 -> plain numeric index 0..10
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "1").Value = $i
}

# Row 2 -> blank, except E2 = "Drive"
$ws.Range("E2").Value = "Drive"

# Row 3 -> old header text (minus thread_size / material_surface)
$ws.Range("A3").Value = $oldHeaders["A"]
$ws.Range("B3").Value = $oldHeaders["B"]
$ws.Range("C3").Value = $oldHeaders["C"]
$ws.Range("D3").Value = $oldHeaders["D"]
$ws.Range("E3").Value = $oldHeaders["E"]
$ws.Range("F3").Value = $oldHeaders["F"]
$ws.Range("G3").Value = $oldHeaders["G"]
$ws.Range("I3").Value = $oldHeaders["I"]

Write-Host "Done restructuring header rows"
